# "added input output to splitter"
#
# On the "Serviços In Out" sheet, fill in the Estado/Input/Output columns
# for the "Split" task (row 5) and the Input column for the "Agregate"
# task (row 6). Also move the active/selected tab from "Tasks" to
# "Serviços In Out".

$wb = $excel.ActiveWorkbook

$wsTasks = $wb.Worksheets.Item("Tasks")
$wsServ  = $wb.Worksheets.Item("Serviços In Out")

# Row 5 - "Split": Estado / Input / Output
$wsServ.Range("B5").Value = "Done"
$wsServ.Range("C5").Value = "Input da webprobe "
$wsServ.Range("D5").Value = "Split da mensagem em companies esb que são passadas assynconamente para o serviço configurado no esb"

# Row 6 - "Agregate": Estado / Input
$wsServ.Range("B6").Value = "On Going"
$wsServ.Range("C6").Value = "Mensagem com as companhias partidas"

# Row 5 grows taller to fit the wrapped Output text
$wsServ.Range("A5:D5").RowHeight = 45
$wsServ.Range("D5").WrapText = $true

# The "Serviços In Out" tab becomes the selected/active one (previously
# "Tasks" was tabSelected), with the cursor resting on C6.
$wsServ.Range("C6").Select()
$wsServ.Activate()
